$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.115.35'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.795.71'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.50%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.07'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5345'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3768'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07459'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.71%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.098'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9989'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.67'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.115'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.18%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.238'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.15%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.782.22'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.12'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001056'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06462'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.38'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9979'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.908'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.132.12'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.03%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.101'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.12'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.26'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.998.18'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.300'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.46'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.117'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1048'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.651'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.569'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2259'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06488'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.46%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.70%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.025'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.500'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.448'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6166'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.08'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.174'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9976'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.28'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.673'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5771'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.34'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.190'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.929'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06811'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.45%  '
